# Add one new week of "Apio" price data at the top of the historical
# series (rows 367-368), pushing the existing data down by two rows.
#
# The sheet is a long historical log of weekly price observations, two
# rows per week (Primera / Segunda quality grades). A new week's worth
# of observations is inserted right before the current row 367, shifting
# all subsequent rows down by 2 (dimension grows from R472 to R474).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 367 (shifts rows 367:472 down to 369:474).
$ws.Rows.Item(367).Insert()
$ws.Rows.Item(367).Insert()

# --- New row 367: "Primera" quality grade for the new week ---
$ws.Cells.Item(367, 1).Value2  = 9
$ws.Cells.Item(367, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(367, 3).Value2  = "Metropolitana"
$ws.Cells.Item(367, 4).Value2  = 45093
$ws.Cells.Item(367, 5).Value2  = 13
$ws.Cells.Item(367, 6).Value2  = 100112017
$ws.Cells.Item(367, 7).Value2  = "Apio"
$ws.Cells.Item(367, 8).Value2  = "Americana (o)"
$ws.Cells.Item(367, 9).Value2  = "Primera"
$ws.Cells.Item(367, 10).Value2 = 70
$ws.Cells.Item(367, 11).Value2 = 7000
$ws.Cells.Item(367, 12).Value2 = 8000
$ws.Cells.Item(367, 13).Value2 = 7500
$ws.Cells.Item(367, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(367, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(367, 16).Value2 = 1250
$ws.Cells.Item(367, 17).Value2 = 6
$ws.Cells.Item(367, 18).Value2 = "Hortaliza"

# --- New row 368: "Segunda" quality grade for the new week ---
$ws.Cells.Item(368, 1).Value2  = 9
$ws.Cells.Item(368, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(368, 3).Value2  = "Metropolitana"
$ws.Cells.Item(368, 4).Value2  = 45093
$ws.Cells.Item(368, 5).Value2  = 13
$ws.Cells.Item(368, 6).Value2  = 100112017
$ws.Cells.Item(368, 7).Value2  = "Apio"
$ws.Cells.Item(368, 8).Value2  = "Americana (o)"
$ws.Cells.Item(368, 9).Value2  = "Segunda"
$ws.Cells.Item(368, 10).Value2 = 52
$ws.Cells.Item(368, 11).Value2 = 6000
$ws.Cells.Item(368, 12).Value2 = 6000
$ws.Cells.Item(368, 13).Value2 = 6000
$ws.Cells.Item(368, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(368, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(368, 16).Value2 = 1000
$ws.Cells.Item(368, 17).Value2 = 6
$ws.Cells.Item(368, 18).Value2 = "Hortaliza"
